$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Weekly crime statistics table updates (rows 14-30) ---

# Row 14
$ws.Range("C14").Value = 3
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 57
$ws.Range("K14").Value = 21.276595744680
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -18.571428571428
$ws.Range("N14").Value = -73.611111111111

# Row 15
$ws.Range("C15").Value = 7
$ws.Range("E15").Value = 40
$ws.Range("F15").Value = 20
$ws.Range("H15").Value = 11.111111111111
$ws.Range("I15").Value = 185
$ws.Range("J15").Value = 187
$ws.Range("K15").Value = -1.069518716577
$ws.Range("L15").Value = -12.322274881516
$ws.Range("M15").Value = 20.129870129870
$ws.Range("N15").Value = -62.090163934426

# Row 16
$ws.Range("C16").Value = 39
$ws.Range("D16").Value = 47
$ws.Range("E16").Value = -17.021276595744
$ws.Range("F16").Value = 151
$ws.Range("G16").Value = 158
$ws.Range("H16").Value = -4.430379746835
$ws.Range("I16").Value = 1665
$ws.Range("J16").Value = 1195
$ws.Range("K16").Value = 39.330543933054
$ws.Range("L16").Value = 19.354838709677
$ws.Range("M16").Value = -32.917002417405
$ws.Range("N16").Value = -86.027190332326

# Row 17
$ws.Range("D17").Value = 62
$ws.Range("E17").Value = -1.612903225806
$ws.Range("F17").Value = 255
$ws.Range("G17").Value = 277
$ws.Range("H17").Value = -7.942238267148
$ws.Range("I17").Value = 2936
$ws.Range("J17").Value = 2657
$ws.Range("K17").Value = 10.500564546481
$ws.Range("L17").Value = 21.674264401160
$ws.Range("M17").Value = 40.009537434430
$ws.Range("N17").Value = -49.256826823366

# Row 18
$ws.Range("C18").Value = 44
$ws.Range("D18").Value = 38
$ws.Range("E18").Value = 15.789473684210
$ws.Range("F18").Value = 187
$ws.Range("G18").Value = 150
$ws.Range("H18").Value = 24.666666666666
$ws.Range("I18").Value = 1821
$ws.Range("J18").Value = 1493
$ws.Range("K18").Value = 21.969189551239
$ws.Range("L18").Value = 8.781362007168
$ws.Range("M18").Value = -38.748738647830
$ws.Range("N18").Value = -88.472494777489

# Row 19
$ws.Range("C19").Value = 124
$ws.Range("E19").Value = -3.875968992248
$ws.Range("F19").Value = 515
$ws.Range("G19").Value = 502
$ws.Range("H19").Value = 2.589641434262
$ws.Range("I19").Value = 6092
$ws.Range("J19").Value = 4315
$ws.Range("K19").Value = 41.181923522595
$ws.Range("L19").Value = 51.166253101737
$ws.Range("M19").Value = 30.254436604661
$ws.Range("N19").Value = -19.471249173826

# Row 20
$ws.Range("C20").Value = 35
$ws.Range("D20").Value = 32
$ws.Range("E20").Value = 9.375
$ws.Range("F20").Value = 146
$ws.Range("G20").Value = 127
$ws.Range("H20").Value = 14.960629921259
$ws.Range("I20").Value = 1564
$ws.Range("J20").Value = 1113
$ws.Range("K20").Value = 40.521114106019
$ws.Range("L20").Value = 27.986906710311
$ws.Range("M20").Value = -5.441354292623
$ws.Range("N20").Value = -92.115743307959

# Row 21
$ws.Range("C21").Value = 313
$ws.Range("D21").Value = 313
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1280
$ws.Range("G21").Value = 1238
$ws.Range("H21").Value = 3.392568659127
$ws.Range("I21").Value = 14320
$ws.Range("J21").Value = 11007
$ws.Range("K21").Value = 30.099027891341
$ws.Range("L21").Value = 30.158153063079
$ws.Range("M21").Value = 1.509888707733
$ws.Range("N21").Value = -76.755133511890

# Row 22
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 4
$ws.Range("F22").Value = 17
$ws.Range("G22").Value = 13
$ws.Range("H22").Value = 30.769230769230
$ws.Range("I22").Value = 163
$ws.Range("J22").Value = 125
$ws.Range("K22").Value = 30.4
$ws.Range("L22").Value = 0.617283950617
$ws.Range("M22").Value = -32.921810699588

# Row 23
$ws.Range("C23").Value = 13
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 160
$ws.Range("F23").Value = 46
$ws.Range("G23").Value = 32
$ws.Range("H23").Value = 43.75
$ws.Range("I23").Value = 477
$ws.Range("J23").Value = 408
$ws.Range("K23").Value = 16.911764705882
$ws.Range("L23").Value = 39.067055393586
$ws.Range("M23").Value = 63.356164383561

# Row 24
$ws.Range("C24").Value = 325
$ws.Range("D24").Value = 258
$ws.Range("E24").Value = 25.968992248062
$ws.Range("F24").Value = 1302
$ws.Range("G24").Value = 1132
$ws.Range("H24").Value = 15.017667844523
$ws.Range("I24").Value = 13887
$ws.Range("J24").Value = 10095
$ws.Range("K24").Value = 37.563150074294
$ws.Range("L24").Value = 34.772903726708
$ws.Range("M24").Value = 30.960015088645

# Row 25
$ws.Range("C25").Value = 96
$ws.Range("D25").Value = 111
$ws.Range("E25").Value = -13.513513513513
$ws.Range("F25").Value = 408
$ws.Range("G25").Value = 415
$ws.Range("H25").Value = -1.686746987951
$ws.Range("I25").Value = 4742
$ws.Range("J25").Value = 4175
$ws.Range("K25").Value = 13.580838323353
$ws.Range("L25").Value = 24.201152435830
$ws.Range("M25").Value = -18.269562219924

# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -18.181818181818
$ws.Range("F26").Value = 33
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = 17.857142857142
$ws.Range("I26").Value = 295
$ws.Range("J26").Value = 293
$ws.Range("K26").Value = 0.682593856655
$ws.Range("L26").Value = -4.530744336569

# Row 27
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 73
$ws.Range("G27").Value = 59
$ws.Range("H27").Value = 23.728813559322
$ws.Range("I27").Value = 608
$ws.Range("J27").Value = 541
$ws.Range("K27").Value = 12.384473197781
$ws.Range("L27").Value = 35.714285714285

# Row 28
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 18
$ws.Range("G28").Value = 19
$ws.Range("H28").Value = -5.263157894736
$ws.Range("I28").Value = 194
$ws.Range("J28").Value = 176
$ws.Range("K28").Value = 10.227272727272
$ws.Range("L28").Value = -25.384615384615
$ws.Range("M28").Value = -17.446808510638
$ws.Range("N28").Value = -71.802325581395

# Row 29
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 66.666666666666
$ws.Range("F29").Value = 17
$ws.Range("G29").Value = 15
$ws.Range("H29").Value = 13.333333333333
$ws.Range("I29").Value = 152
$ws.Range("J29").Value = 154
$ws.Range("K29").Value = -1.298701298701
$ws.Range("L29").Value = -22.842639593908
$ws.Range("M29").Value = -22.051282051282
$ws.Range("N29").Value = -74.708818635607

# Row 30
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 4
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 6
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = -45.454545454545
$ws.Range("I30").Value = 100
$ws.Range("J30").Value = 63
$ws.Range("K30").Value = 58.730158730158
$ws.Range("L30").Value = 138.095238095238
